$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 12849.5 pesos`n✅ 12849.5 pesos = 3.33 = 976.3 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 299
$ws2.Range("O10").Value = 3842
$ws2.Range("N12").Value = 3855
$ws2.Range("O12").Value = 292.9
